$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price-report row is inserted at row 166 (pushing the existing
# rows 166-186 down to 167-187), bringing the sheet from A1:R186 to A1:R187.
$ws.Rows.Item(166).Insert()

# Populate the newly inserted row with the latest "Berenjena" observation.
$ws.Cells.Item(166, 1).Value = 10
$ws.Cells.Item(166, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(166, 3).Value = "La Araucanía"
$ws.Cells.Item(166, 4).Value = 44474
$ws.Cells.Item(166, 5).Value = 9
$ws.Cells.Item(166, 6).Value = 100112001
$ws.Cells.Item(166, 7).Value = "Berenjena"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 30
$ws.Cells.Item(166, 11).Value = 10000
$ws.Cells.Item(166, 12).Value = 10000
$ws.Cells.Item(166, 13).Value = 10000
$ws.Cells.Item(166, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(166, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(166, 16).Value = 167
$ws.Cells.Item(166, 17).Value = 60
$ws.Cells.Item(166, 18).Value = "Hortaliza"
